$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append (date serial, B, C, D)
$newRows = @(
    @(358, 44432, 0, 2, 47.65308553728854),
    @(359, 44433, 0, 2, 47.65308553728854),
    @(360, 44434, 0, 2, 47.65308553728854),
    @(361, 44435, 1, 2, 47.65308553728854),
    @(362, 44436, 0, 2, 47.65308553728854),
    @(363, 44437, 0, 2, 47.65308553728854),
    @(364, 44438, 1, 2, 47.65308553728854),
    @(365, 44439, 0, 2, 47.65308553728854),
    @(366, 44440, 0, 2, 47.65308553728854)
)

foreach ($r in $newRows) {
    $rowIdx = $r[0]
    $dateVal = $r[1]
    $bVal = $r[2]
    $cVal = $r[3]
    $dVal = $r[4]

    # Copy formatting from the last existing data row (357) into the new row
    $ws.Range("A357:D357").Copy() | Out-Null
    $destRow = $ws.Range("A" + $rowIdx + ":D" + $rowIdx)
    $destRow.PasteSpecial(-4122) | Out-Null

    $ws.Cells.Item($rowIdx, 1).Value2 = $dateVal
    $ws.Cells.Item($rowIdx, 2).Value2 = $bVal
    $ws.Cells.Item($rowIdx, 3).Value2 = $cVal
    $ws.Cells.Item($rowIdx, 4).Value2 = $dVal
}

$excel.CutCopyMode = 0
